$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1, shared string) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 18:20"

# --- Update country case numbers that changed in place (country stays the same) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 191193
$ws.Range("C4").Value = 2663
$ws.Range("D4").Value = 7531
$ws.Range("E4").Value = 179524
$ws.Range("F4").Value = 4576
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = 4138

# Row 5: Italia
$ws.Range("B5").Value = 110574
$ws.Range("C5").Value = 4782
$ws.Range("D5").Value = 16847
$ws.Range("E5").Value = 80572
$ws.Range("F5").Value = 4035
$ws.Range("G5").Value = 727
$ws.Range("H5").Value = 13155

# Row 8: Alemania
$ws.Range("B8").Value = 75754
$ws.Range("C8").Value = 3946
$ws.Range("D8").Value = 18700
$ws.Range("E8").Value = 56206
$ws.Range("F8").Value = 3405
$ws.Range("G8").Value = 73
$ws.Range("H8").Value = 848

# Row 16: Austria
$ws.Range("B16").Value = 10585
$ws.Range("C16").Value = 405
$ws.Range("E16").Value = 9003

# Row 18: Canada
$ws.Range("B18").Value = 9017
$ws.Range("C18").Value = 405
$ws.Range("D18").Value = 1445
$ws.Range("E18").Value = 7464
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 108

# Row 20: Brasil
$ws.Range("B20").Value = 5916
$ws.Range("C20").Value = 199
$ws.Range("E20").Value = 5583
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 206

# Row 25: Chequia
$ws.Range("D25").Value = 61
$ws.Range("E25").Value = 3408

# Row 51: Argentina
$ws.Range("E51").Value = 778
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 28

# --- Rows 55-59: countries shift up one slot as Argelia moves up the ranking ---
# Row 55 becomes Argelia (was Eslovenia)
$ws.Range("A55").Value = "Argelia"
$ws.Range("B55").Value = 846
$ws.Range("C55").Value = 130
$ws.Range("D55").Value = 46
$ws.Range("E55").Value = 742
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 14
$ws.Range("H55").Value = 58

# Row 56 becomes Eslovenia (was Catar)
$ws.Range("A56").Value = "Eslovenia"
$ws.Range("B56").Value = 841
$ws.Range("C56").Value = 39
$ws.Range("D56").Value = 10
$ws.Range("E56").Value = 816
$ws.Range("F56").Value = 31
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 15

# Row 57 becomes Catar (was Estonia)
$ws.Range("A57").Value = "Catar"
$ws.Range("B57").Value = 781
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 62
$ws.Range("E57").Value = 717
$ws.Range("F57").Value = 6
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 2

# Row 58 becomes Estonia (was Hong Kong)
$ws.Range("A58").Value = "Estonia"
$ws.Range("B58").Value = 779
$ws.Range("C58").Value = 34
$ws.Range("D58").Value = 33
$ws.Range("E58").Value = 741
$ws.Range("F58").Value = 15
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 5

# Row 59 becomes Hong Kong (was Argelia)
$ws.Range("A59").Value = "Hong Kong"
$ws.Range("B59").Value = 765
$ws.Range("C59").Value = 50
$ws.Range("D59").Value = 147
$ws.Range("E59").Value = 614
$ws.Range("F59").Value = 5
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 4

# Row 85: Republica de Chipre
$ws.Range("F85").Value = 11
